# Update the cryptocurrency price/volume table with freshly scraped values.
# Columns D (Price) and E (Volume(1h)) hold text-formatted figures (e.g.
# "226.80", "  +2.50%  ") that must stay as text rather than being coerced
# into numbers/percentages by Excel, so we force Text number formatting
# while assigning the values, then restore the original cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$dataRange = $ws.Range('D2:E51')
$origStyle = $dataRange.Style
$dataRange.NumberFormat = '@'

$ws.Range('D2').Value = '39.740.34'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '2.161.36'
$ws.Range('E3').Value = '  +2.92%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '226.80'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = '0.622'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('D7').Value = '62.89'
$ws.Range('E7').Value = '  +1.62%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').Value = '0.0840'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '15.88'
$ws.Range('E12').Value = '  +0.74%  '
$ws.Range('D13').Value = '2.481.08'
$ws.Range('E13').Value = '  +2.88%  '
$ws.Range('D14').Value = '21.80'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '5.51'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '2.165.94'
$ws.Range('E17').Value = '  +2.78%  '
$ws.Range('D18').Value = '39.740.54'
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('D19').Value = '71.77'
$ws.Range('E19').Value = '  +0.23%  '
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').Value = '228.19'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '2.36'
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('D26').Value = '172.18'
$ws.Range('E26').Value = '  +1.17%  '
$ws.Range('D27').Value = '9.43'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('E28').Value = '  +2.16%  '
$ws.Range('D29').Value = '1.44'
$ws.Range('E29').Value = '  +2.11%  '
$ws.Range('D30').Value = '19.59'
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('E31').Value = '  +5.82%  '
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  +0.67%  '
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('E34').Value = '  -2.17%  '
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('D36').Value = '0.0616'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('E37').Value = '  +7.50%  '
$ws.Range('E38').Value = '  +1.73%  '
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').Value = '4.91'
$ws.Range('E40').Value = '  +17.82%  '
$ws.Range('D41').Value = '102.75'
$ws.Range('E41').Value = '  +1.20%  '
$ws.Range('D42').Value = '0.0227'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('D43').Value = '17.56'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').Value = '1.511.90'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('E45').Value = '  -0.09%  '
$ws.Range('D46').Value = '7.93'
$ws.Range('E46').Value = '  +2.08%  '
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').Value = '2.81'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  +0.88%  '
$ws.Range('D50').Value = '49.61'
$ws.Range('E50').Value = '  +8.09%  '
$ws.Range('D51').Value = '0.000190'
$ws.Range('E51').Value = '  +28.02%  '

$dataRange.Style = $origStyle
